# Apply updated crypto price / volume figures (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D values that look numeric must stay text, matching the sheet's
# original inline-string storage (e.g. '1.00' must not become 1).
$textCells = @("D4", "D5", "D6", "D10", "D14", "D17", "D19", "D20", "D21", "D22", "D23", "D24", "D27", "D28", "D29", "D31", "D34", "D36", "D39", "D40", "D43", "D44", "D45", "D49", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '63.323.38'
$ws.Range("E2").Value = '  +1.74%  '
$ws.Range("D3").Value = '3.472.80'
$ws.Range("E3").Value = '  +1.25%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '581.96'
$ws.Range("E5").Value = '  +0.35%  '
$ws.Range("D6").Value = '147.67'
$ws.Range("E6").Value = '  +1.34%  '
$ws.Range("D7").Value = '3.472.91'
$ws.Range("E7").Value = '  +1.23%  '
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("E9").Value = '  +0.46%  '
$ws.Range("D10").Value = '7.73'
$ws.Range("E10").Value = '  +1.62%  '
$ws.Range("E11").Value = '  +0.74%  '
$ws.Range("E12").Value = '  +4.52%  '
$ws.Range("D13").Value = '4.068.12'
$ws.Range("E13").Value = '  +1.31%  '
$ws.Range("D14").Value = '29.52'
$ws.Range("E14").Value = '  +2.22%  '
$ws.Range("E15").Value = '  +2.62%  '
$ws.Range("D16").Value = '3.467.01'
$ws.Range("E16").Value = '  +1.13%  '
$ws.Range("D17").Value = '0.0000172'
$ws.Range("E17").Value = '  +0.64%  '
$ws.Range("D18").Value = '63.351.83'
$ws.Range("E18").Value = '  +1.82%  '
$ws.Range("D19").Value = '6.39'
$ws.Range("E19").Value = '  +3.01%  '
$ws.Range("D20").Value = '14.50'
$ws.Range("E20").Value = '  +2.93%  '
$ws.Range("D21").Value = '9.36'
$ws.Range("E21").Value = '  +1.35%  '
$ws.Range("D22").Value = '389.21'
$ws.Range("E22").Value = '  -1.32%  '
$ws.Range("D23").Value = '0.566'
$ws.Range("E23").Value = '  +1.91%  '
$ws.Range("D24").Value = '74.56'
$ws.Range("E24").Value = '  -0.48%  '
$ws.Range("E25").Value = '  -0.08%  '
$ws.Range("D26").Value = '3.615.03'
$ws.Range("E26").Value = '  +1.35%  '
$ws.Range("D27").Value = '0.0000116'
$ws.Range("E27").Value = '  -0.33%  '
$ws.Range("D28").Value = '0.182'
$ws.Range("E28").Value = '  -1.98%  '
$ws.Range("D29").Value = '7.65'
$ws.Range("E29").Value = '  +1.31%  '
$ws.Range("E30").Value = '  +0.01%  '
$ws.Range("D31").Value = '8.22'
$ws.Range("E31").Value = '  +2.13%  '
$ws.Range("E32").Value = '  -0.31%  '
$ws.Range("E33").Value = '  +0.07%  '
$ws.Range("D34").Value = '1.35'
$ws.Range("E34").Value = '  -4.82%  '
$ws.Range("E35").Value = '  -0.66%  '
$ws.Range("D36").Value = '5.34'
$ws.Range("E36").Value = '  +0.84%  '
$ws.Range("E37").Value = '  +2.03%  '
$ws.Range("E38").Value = '  +6.31%  '
$ws.Range("D39").Value = '31.91'
$ws.Range("E39").Value = '  +11.46%  '
$ws.Range("D40").Value = '168.72'
$ws.Range("E40").Value = '  +0.61%  '
$ws.Range("D41").Value = '3.512.08'
$ws.Range("E41").Value = '  +1.44%  '
$ws.Range("E42").Value = '  +1.52%  '
$ws.Range("D43").Value = '0.796'
$ws.Range("E43").Value = '  +0.82%  '
$ws.Range("D44").Value = '1.75'
$ws.Range("E44").Value = '  +3.64%  '
$ws.Range("D45").Value = '42.41'
$ws.Range("E45").Value = '  -1.19%  '
$ws.Range("E46").Value = '  +2.94%  '
$ws.Range("E47").Value = '  -1.31%  '
$ws.Range("D48").Value = '2.597.54'
$ws.Range("E48").Value = '  +3.32%  '
$ws.Range("D49").Value = '2.32'
$ws.Range("E49").Value = '  +8.26%  '
$ws.Range("E50").Value = '  +2.26%  '
$ws.Range("D51").Value = '23.07'
$ws.Range("E51").Value = '  -0.37%  '
